$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "B2" used to hold a stray pandas artifact label ("unnamed: 1_level_1");
# correct it to the real header value "total".
$ws.Range("B2").Value = "total"

# Remove the two empty sub-header rows that had no data under them
# ("situação do domicílio" at row 5 and "grandes regiões e unidades da
# federação" at row 8). Delete the lower one first so row numbers above
# it stay valid while we still need them.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
